$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5603.5
$ws.Range("I6").Value = 5603.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 16810.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -16698.5
$ws.Range("N6").Value = $null

$ws.Range("H33").Value = 202.21053
$ws.Range("I33").Value = 158.28572
$ws.Range("J33").Value = 325.2
$ws.Range("K33").Value = 158.28572
$ws.Range("L33").Value = 325.2
$ws.Range("M33").Value = 70.71428
$ws.Range("N33").Value = -783.2

$ws.Range("H98").Value = 658963.25
$ws.Range("I98").Value = 1015652.4
$ws.Range("J98").Value = 5033.1665
$ws.Range("K98").Value = 1015652.4
$ws.Range("L98").Value = 5033.1665
$ws.Range("M98").Value = -1014154.4
$ws.Range("N98").Value = -8029.1665

$ws.Range("H122").Value = 658963.25
$ws.Range("I122").Value = 1015652.4
$ws.Range("J122").Value = 5033.1665
$ws.Range("K122").Value = 3046957.2
$ws.Range("L122").Value = 15099.4995
$ws.Range("M122").Value = -3044507.2
$ws.Range("N122").Value = -19999.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1190.6666
$ws.Range("I14").Value = 1190.6666
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1190.6666
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1015.6666
$ws.Range("N14").Value = $null

$ws.Range("H32").Value = 2182.463
$ws.Range("I32").Value = 1747.826
$ws.Range("J32").Value = 4681.625
$ws.Range("K32").Value = 1747.826
$ws.Range("L32").Value = 4681.625
$ws.Range("M32").Value = -1460.826
$ws.Range("N32").Value = -5255.625

$ws.Range("H45").Value = 2558
$ws.Range("I45").Value = 1379.1111
$ws.Range("K45").Value = 1379.1111
$ws.Range("M45").Value = -1002.1111

$ws.Range("H74").Value = 821.0345
$ws.Range("I74").Value = 938.94116
$ws.Range("J74").Value = 654
$ws.Range("K74").Value = 938.94116
$ws.Range("L74").Value = 654
$ws.Range("M74").Value = -64.94115999999997
$ws.Range("N74").Value = -2402

$ws.Range("H77").Value = 821.0345
$ws.Range("I77").Value = 938.94116
$ws.Range("J77").Value = 654
$ws.Range("K77").Value = 4694.7058
$ws.Range("L77").Value = 3270
$ws.Range("M77").Value = -326.7057999999997
$ws.Range("N77").Value = -12006

$ws.Range("H132").Value = 4002.15
$ws.Range("I132").Value = 3724
$ws.Range("K132").Value = 11172
$ws.Range("M132").Value = -8642

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1850.6923
$ws.Range("I31").Value = 1082.7059
$ws.Range("J31").Value = 2444.1365
$ws.Range("K31").Value = 1082.7059
$ws.Range("L31").Value = 2444.1365
$ws.Range("M31").Value = -787.7058999999999
$ws.Range("N31").Value = -3034.1365

$ws.Range("H34").Value = 1850.6923
$ws.Range("I34").Value = 1082.7059
$ws.Range("J34").Value = 2444.1365
$ws.Range("K34").Value = 1082.7059
$ws.Range("L34").Value = 2444.1365
$ws.Range("M34").Value = -880.7058999999999
$ws.Range("N34").Value = -2848.1365

$ws.Range("H105").Value = 501.54544
$ws.Range("I105").Value = 464.625
$ws.Range("K105").Value = 464.625
$ws.Range("M105").Value = 1282.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1066.0769
$ws.Range("I11").Value = 341.72726
$ws.Range("K11").Value = 1025.18178
$ws.Range("M11").Value = -885.1817799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2599900
$ws.Range("I14").Value = 2599900
$ws.Range("K14").Value = 2599900
$ws.Range("M14").Value = -2599732

$ws.Range("H80").Value = 2350
$ws.Range("I80").Value = 2364
$ws.Range("J80").Value = 2262.5
$ws.Range("K80").Value = 2364
$ws.Range("L80").Value = 2262.5
$ws.Range("M80").Value = -1366
$ws.Range("N80").Value = -4258.5

$ws.Range("H83").Value = 2350
$ws.Range("I83").Value = 2364
$ws.Range("J83").Value = 2262.5
$ws.Range("K83").Value = 11820
$ws.Range("L83").Value = 11312.5
$ws.Range("M83").Value = -6828
$ws.Range("N83").Value = -21296.5

$ws.Range("H97").Value = 2771.9285
$ws.Range("I97").Value = 2368.923
$ws.Range("J97").Value = 8011
$ws.Range("K97").Value = 2368.923
$ws.Range("L97").Value = 8011
$ws.Range("M97").Value = -1872.923
$ws.Range("N97").Value = -9003

$ws.Range("H107").Value = 1006
$ws.Range("I107").Value = 1006
$ws.Range("K107").Value = 1006
$ws.Range("M107").Value = 914

$ws.Range("H122").Value = 2225903.5
$ws.Range("I122").Value = 3707672.8
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 11123018.4
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -11120568.4
$ws.Range("N122").Value = -14650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3194.7932
$ws.Range("J7").Value = 3441.1304
$ws.Range("L7").Value = 3441.1304
$ws.Range("N7").Value = -3665.1304

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null

$ws.Range("H68").Value = 2320
$ws.Range("I68").Value = 1933.3334
$ws.Range("J68").Value = 2900
$ws.Range("K68").Value = 1933.3334
$ws.Range("L68").Value = 2900
$ws.Range("M68").Value = -1184.3334
$ws.Range("N68").Value = -4398

$ws.Range("H71").Value = 2320
$ws.Range("I71").Value = 1933.3334
$ws.Range("J71").Value = 2900
$ws.Range("K71").Value = 9666.666999999999
$ws.Range("L71").Value = 14500
$ws.Range("M71").Value = -5922.666999999999
$ws.Range("N71").Value = -21988

$ws.Range("H82").Value = 46861.637
$ws.Range("J82").Value = 1745
$ws.Range("L82").Value = 1745
$ws.Range("N82").Value = -2467

$ws.Range("H85").Value = 46861.637
$ws.Range("J85").Value = 1745
$ws.Range("L85").Value = 1745
$ws.Range("N85").Value = -4241

$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws.Range("H126").Value = 3194.7932
$ws.Range("J126").Value = 3441.1304
$ws.Range("L126").Value = 10323.3912
$ws.Range("N126").Value = -15263.3912

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 18666.666
$ws.Range("I7").Value = 50000
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 50000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -49887
$ws.Range("N7").Value = -3226

$ws.Range("H62").Value = 31512.285
$ws.Range("I62").Value = 100000
$ws.Range("K62").Value = 100000
$ws.Range("M62").Value = -99376

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = $null

$ws.Range("H65").Value = 31512.285
$ws.Range("I65").Value = 100000
$ws.Range("K65").Value = 500000
$ws.Range("M65").Value = -496880

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = $null

$ws.Range("H121").Value = 26033.846
$ws.Range("J121").Value = 26033.846
$ws.Range("L121").Value = 26033.846
$ws.Range("N121").Value = -29527.846

$ws.Range("H132").Value = 18520952
$ws.Range("I132").Value = 23811288
$ws.Range("J132").Value = 4771.5
$ws.Range("K132").Value = 71433864
$ws.Range("L132").Value = 14314.5
$ws.Range("M132").Value = -71431334
$ws.Range("N132").Value = -19374.5
